$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Enter the new values for rows 39-47 (new TESTING section)
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "TESTING"

$ws.Range("A41").Value = "Updated Task Requirements and Project Requirements"
$ws.Range("B41").Value = "Will Maberry"
$ws.Range("C41").Value = "Will Maberry"

$ws.Range("A42").Value = "bug fixing game functionality and game user story test"
$ws.Range("B42").Value = "Will Maberry"
$ws.Range("C42").Value = "Will Maberry"

$ws.Range("A43").Value = "bug fixing online players and online players user story test"
$ws.Range("B43").Value = "Grace Daily"
$ws.Range("C43").Value = "Grace Daily"

$ws.Range("A44").Value = "bug fixing message sending and message user story test"
$ws.Range("B44").Value = "David Oyekola"
$ws.Range("C44").Value = "David Oyekola"

$ws.Range("A45").Value = "bug fixing leaderboard functionality and leaderboard user story test"
$ws.Range("B45").Value = "Subodh Neupane"
$ws.Range("C45").Value = "Will Maberry coded, wrote automatic tests, and user story test"

$ws.Range("A46").Value = "bug fixing user authentication and user authentication story test"
$ws.Range("B46").Value = "Cody Mercer"
$ws.Range("C46").Value = "Cody Mercer"

$ws.Range("A47").Value = "Automatic WholeGameTest.java"
$ws.Range("B47").Value = "Ammar Rafiq"
$ws.Range("C47").Value = "Ammar Rafiq"

# ---------------------------------------------------------------------------
# 2. Copy formatting from the analogous rows in the IMPLEMENTATION section so
#    the new TESTING section is styled the same way (banded rows + header).
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
  $ws.Range($srcAddr).Copy()
  $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# blank separator row (matches row 21)
Copy-Format "A21" "A39"
Copy-Format "B21" "B39"
Copy-Format "C21" "C39"

# section header row (matches row 22)
Copy-Format "A22" "A40"
Copy-Format "B22" "B40"
Copy-Format "C22" "C40"

# data rows (matches rows 23-29 banding pattern)
Copy-Format "A23" "A41"
Copy-Format "B23" "B41"
Copy-Format "C23" "C41"

Copy-Format "A24" "A42"
Copy-Format "B24" "B42"
Copy-Format "C24" "C42"

Copy-Format "A25" "A43"
Copy-Format "B25" "B43"
Copy-Format "C25" "C43"

Copy-Format "A26" "A44"
Copy-Format "B26" "B44"
Copy-Format "C26" "C44"

Copy-Format "A27" "A45"
Copy-Format "B27" "B45"
Copy-Format "C27" "C45"

Copy-Format "A28" "A46"
Copy-Format "B28" "B46"
Copy-Format "C28" "C46"

Copy-Format "A29" "A47"
Copy-Format "B29" "B47"
Copy-Format "C29" "C47"

# ---------------------------------------------------------------------------
# 3. Column width tweaks (A wider to fit new text, C wider & bestFit)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 58
$ws.Columns("C").ColumnWidth = 55.140625

# ---------------------------------------------------------------------------
# 4. Restore clean (non-clipboard) selection / view state
# ---------------------------------------------------------------------------
$excel.CutCopyMode = 0
$ws.Range("C50").Select()
